$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - Worksheets index 1 / sheet1.xml
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 52
$ws1.Range("F4").Value = 1582
$ws1.Range("F5").Value = 280
$ws1.Range("F6").Value = 65
$ws1.Range("F7").Value = 1739
$ws1.Range("F8").Value = 10242
$ws1.Range("F9").Value = 175
$ws1.Range("F10").Value = 139
$ws1.Range("F11").Value = 259
$ws1.Range("F12").Value = 195
$ws1.Range("F13").Value = 393
$ws1.Range("F14").Value = 7102
$ws1.Range("F15").Value = 1104
$ws1.Range("F16").Value = 666
$ws1.Range("F17").Value = 55
$ws1.Range("F19").Value = 241

# Sheet "演出" (Performances) - Worksheets index 2 / sheet2.xml
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 10
$ws2.Range("F3").Value = 557

# Sheet "全部类型" (All types) - Worksheets index 4 / sheet4.xml
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 52
$ws4.Range("F4").Value = 1582
$ws4.Range("F5").Value = 280
$ws4.Range("F6").Value = 10
$ws4.Range("F7").Value = 65
$ws4.Range("F8").Value = 1739
$ws4.Range("F9").Value = 557
$ws4.Range("F11").Value = 10242
$ws4.Range("F12").Value = 175
$ws4.Range("F13").Value = 139
$ws4.Range("F14").Value = 259
$ws4.Range("F15").Value = 195
$ws4.Range("F16").Value = 393
$ws4.Range("F17").Value = 7102
$ws4.Range("F18").Value = 1104
$ws4.Range("F19").Value = 666
$ws4.Range("F20").Value = 55
$ws4.Range("F22").Value = 241
